# Restored from revision of admin on 03/06/2020 at 01:50:18 PM.TEST Author: admin. Type: SAVE.
# Sets cell C10 on the active ("Rules") sheet from 18 to 1 (keeps existing
# number formatting / style; only the stored value changes).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1.0
